$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Insert a new row at position 4 (pushes old row4 "disability
#    persons" row down to row5, and old row5 "source" row to row6)
# ---------------------------------------------------------------
$ws.Rows.Item(4).Insert()

# ---------------------------------------------------------------
# 2. Update text content
# ---------------------------------------------------------------
# Row2 unchanged: "(End of year, persons)"

# New row4: "family with disabilities Persons "
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("B4").Value = 1492
$ws.Range("C4").Value = 1494
$ws.Range("D4").Value = 1500
$ws.Range("E4").Value = 1585
$ws.Range("F4").Value = 1636
$ws.Range("G4").Value = 1651
$ws.Range("H4").Value = 1601
$ws.Range("I4").Value = 1618

# Row5 (previously row4): rename text, update values
$ws.Range("A5").Value = "disabilities Persons "
$ws.Range("B5").Value = 1651
$ws.Range("C5").Value = 1658
$ws.Range("D5").Value = 1647
$ws.Range("E5").Value = 1742
$ws.Range("F5").Value = 1787
$ws.Range("G5").Value = 1794
$ws.Range("H5").Value = 1738
$ws.Range("I5").Value = 1757

# Row6 (previously row5/source row) keeps its original "Source: ..." rich
# text - no change needed there.

# Title row1 gets the new title text (was: "Number of Internally Displaced
# Disability Persons Receiving Social Package in Zugdidi Municipality")
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Zugdidi Municipality"

Write-Host "content updated"

# ---------------------------------------------------------------
# 3. Merge the title cells A1:I1
# ---------------------------------------------------------------
$ws.Range("A1:I1").Merge()

# ---------------------------------------------------------------
# 4. Row heights
# ---------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 51
$ws.Rows.Item(2).RowHeight = $ws.StandardHeight
$ws.Rows.Item(4).RowHeight = 24.75
$ws.Rows.Item(6).RowHeight = 27.75

# ---------------------------------------------------------------
# 5. Cell formatting
# ---------------------------------------------------------------

# Row1 title: centre horizontally too, and wrap
$ws.Range("A1:I1").HorizontalAlignment = -4108
$ws.Range("A1:I1").VerticalAlignment = -4108
$ws.Range("A1:I1").WrapText = $true

# Row3 A3: font becomes Sylfaen 11 (was Arial 10)
$ws.Range("A3").Font.Name = "Sylfaen"
$ws.Range("A3").Font.Size = 11

# Row4 (new): A4 label formatting - fill, bottom border, left/center/wrap
$ws.Range("A4").Interior.Pattern = 1
$ws.Range("A4").Interior.ThemeColor = 0
$ws.Range("A4").Interior.TintAndShade = 0
$ws.Range("A4").Borders.Item(9).LineStyle = 1
$ws.Range("A4").Borders.Item(9).Weight = 2
$ws.Range("A4").Borders.Item(9).ColorIndex = 64
$ws.Range("A4").Borders.Item(8).LineStyle = -4142
$ws.Range("A4").HorizontalAlignment = -4131
$ws.Range("A4").VerticalAlignment = -4108
$ws.Range("A4").WrapText = $true

# Row4 values: plain numeric style (no border / no special alignment)
$ws.Range("B4:I4").NumberFormat = "#\ ##0"
$ws.Range("B4:I4").HorizontalAlignment = -4131
$ws.Range("B4:I4").WrapText = $false
$ws.Range("B4:I4").Interior.Pattern = 1
$ws.Range("B4:I4").Interior.ThemeColor = 0
$ws.Range("B4:I4").Interior.TintAndShade = 0
$ws.Range("B4:I4").Borders.Item(8).LineStyle = -4142
$ws.Range("B4:I4").Borders.Item(9).LineStyle = -4142

# Row5 (was row4): A5 - remove top border, keep bottom border only
$ws.Range("A5").Borders.Item(8).LineStyle = -4142
$ws.Range("A5").Borders.Item(9).LineStyle = 1
$ws.Range("A5").Borders.Item(9).Weight = 2
$ws.Range("A5").Borders.Item(9).ColorIndex = 64

# Row5 values: unify style, drop the "align right" + stray borders
$ws.Range("B5:I5").HorizontalAlignment = -4131
$ws.Range("B5:H5").Borders.Item(9).LineStyle = -4142
$ws.Range("I5").Borders.Item(9).LineStyle = 1
$ws.Range("I5").Borders.Item(9).Weight = 2
$ws.Range("I5").Borders.Item(9).ColorIndex = 64

# Row6 (was row5/source row): font changes to 9pt Arial (not bold) black
$ws.Range("A6:H6").Font.Name = "Arial"
$ws.Range("A6:H6").Font.Size = 9
$ws.Range("A6:H6").Font.Bold = $false
$ws.Range("A6:H6").Font.ColorIndex = 1
$ws.Range("A6").Borders.Item(8).LineStyle = -4142
